$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "X_Exper" column (C) entirely - shifts nothing left of it,
# just drops the column and its header/values.
$ws.Range("C1:C11").EntireColumn.Delete()

# Merge the old "Time" / "X" header text into a single "TimeX" label in A1,
# and clear out the now-unused B1 header cell.
$ws.Range("A1").Value = "TimeX"
$ws.Range("B1").ClearContents()
